$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Author cell
$ws.Range("A1").Value = "Laura Braak - T00198405"

# Example -> real content for "Item Name"
$ws.Range("B3").Value = "Items"

# Example -> real content for "Describe Role in game"
$ws.Range("B4").Value = "The Charakter picks up the Items to earn score. "

# Row 7: was "Eg Turn Left" example row, now the first real Internal Functionality entry
$ws.Range("A7").Value = "Default Constructer"
$ws.Range("B7").Value = "Creating Default Item with a default Value "

# Row 12: was "Eg Push" example row, now a real External Outgoing entry
$ws.Range("A12").Value = "Collecting Animation "
$ws.Range("B12").Value = "When collusion with penguin, start a animation "

# Row 16 ("ShouldTurnLeft") is removed entirely
$ws.Range("A16").Value = ""

# Update the remembered selection
[void]$ws.Range("A13").Select()
